# The dataset sheet had a styled header row (em1, em2, em3, Modifiability,
# Maintainability, SQC, repository, version) in row 1, with bold text,
# thin borders and centered alignment. The commit replaces the generated
# CSV/XLSX export with a version that has no header row at all: every data
# row moves up by one, the sheet shrinks from A1:H84 to A1:H83, and none of
# the remaining cells keep the old header formatting (no bold font, no
# border, no centered alignment).
#
# Deleting row 1 outright reproduces exactly that: Excel shifts all the
# rows below it up by one and drops the formatting that lived only on the
# deleted row, while leaving the values/precision of every other cell
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the header row's distinctive formatting (bold font, thin box
# border, centered alignment) is gone before/with the deletion.
$headerRow = $ws.Rows(1)
$headerRow.ClearFormats()

# Remove the header row entirely; everything below shifts up one row.
$headerRow.Delete()

$wb.Save()
